$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Bulk-add bug fix: a recipient row whose "email" is actually not a valid
# email address ("222") is now detected instead of being silently treated
# like the well-formed rows above it. Record the raw numeric value that was
# uploaded for both the id (A3) and the bogus "email" (B3).
$target = $ws.Range("B3")
$ws.Range("A3").Value = 2
$target.Value = 222

# The cell used to be hyperlinked text reading "2@2" (a valid address). Now
# that the cell holds the flagged numeric value 222, keep the hyperlink
# (still pointing at mailto:2@2) but surface the original "2@2" text as the
# hyperlink's display text, so the Management GUI can show the invalid
# value that was actually uploaded.
for ($i = 1; $i -le $ws.Hyperlinks.Count; $i++) {
    $link = $ws.Hyperlinks.Item($i)
    if ($link.Range.Address() -eq $target.Address()) {
        $link.TextToDisplay = "2@2"
        break
    }
}

# Restore the numeric value; updating TextToDisplay re-writes the cell's
# text to match the hyperlink, so re-apply the flagged value afterwards.
$target.Value = 222

# Mirror the selection recorded for this workbook.
$ws.Range("E17").Select()
